$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.857.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.46%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.484.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.10%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.83%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.72%  '

# Row 7
$ws.Range("E7").Value = '  +0.61%  '

# Row 8
$ws.Range("E8").Value = '  -2.55%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.487.33'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.28%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0997'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.02%  '

# Row 11
$ws.Range("E11").Value = '  -0.06%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.11%  '

# Row 13
$ws.Range("E13").Value = '  -4.79%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.935.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.74%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.747.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.53%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.18%  '

# Row 17
$ws.Range("E17").Value = '  -2.95%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.482.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.26%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.49%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.04%  '

# Row 22
$ws.Range("E22").Value = '  -0.09%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.76%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.35%  '

# Row 25
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.71%  '

# Row 26
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.410'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.93%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.995'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.69%  '

# Row 29
$ws.Range("E29").Value = '  -4.26%  '

# Row 30
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.73%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.40%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.25%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.36'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.78%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.65%  '

# Row 37
$ws.Range("E37").Value = '  -10.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.54'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.42%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.59'
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.790'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.52%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.59%  '

# Row 42
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '275.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.90%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.87%  '

# Row 44
$ws.Range("E44").Value = '  +0.41%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.591'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.10%  '

# Row 46
$ws.Range("E46").Value = '  -1.02%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.52%  '

# Row 48
$ws.Range("E48").Value = '  -3.88%  '

# Row 49
$ws.Range("E49").Value = '  -4.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.761.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.92%  '
